$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, bypassing Excel's
# automatic "looks like a number" detection (needed for the "Documento"
# column, which holds digit-only ids that must stay text) while keeping
# the cell's default style - touching NumberFormat directly would mint a
# brand-new style index, which the source file never had. We enter a
# text-literal formula that evaluates to the desired string, then freeze
# the formula result into a plain value via Copy + Paste Values.
function Set-TextValue($rangeAddr, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($rangeAddr).Formula = '="' + $escaped + '"'
    $ws.Range($rangeAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 2 (existing attendance record): the check-in time and resulting
# status were corrected, and the "Curso" value was cleared out.
$ws.Range("G2").ClearContents()
Set-TextValue "H2" "10:52 a. m."
Set-TextValue "I2" "A tiempo"

# --- Row 3 (new attendance record).
Set-TextValue "A3" "Victor Manuel"
Set-TextValue "B3" "Bonilla Gutierrez"
Set-TextValue "C3" "Permiso por Protección Temporal"
Set-TextValue "D3" "4073477"
Set-TextValue "E3" "ADSO"
Set-TextValue "F3" "Tecnologo"
$ws.Range("G3").ClearContents()
Set-TextValue "H3" "10:52 a. m."
Set-TextValue "I3" "A tiempo"

# Extend the conditional formatting on the "Estado" column so it also
# covers the newly added row, keeping the original rules/dxf mapping.
$rules = $ws.Range("I2").FormatConditions
for ($i = 1; $i -le $rules.Count; $i++) {
    $rules.Item($i).ModifyAppliesToRange($ws.Range("I2:I3"))
}
